$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("L1").Value = "hzj-混合调节_20170516_152754_ASIC_EEG"
$ws.Range("M1").Value = "hzj-混合调节_20170518_134207_ASIC_EEG"
$ws.Range("N1").Value = "hzj-混合调节_20170519_135415_ASIC_EEG"
$ws.Range("O1").Value = "zyx-混合调节_20170516_111228_ASIC_EEG"
$ws.Range("P1").Value = "zyx-混合调节_20170517_110944_ASIC_EEG"
$ws.Range("Q1").Value = "zyx-混合调节_20170518_112337_ASIC_EEG"
$ws.Range("R1").Value = "zyx-混合调节_20170519_124954_ASIC_EEG"
$ws.Range("S1").Value = "zyx-混合调节_20170522_111557_ASIC_EEG"

$ws.Range("L2").Value = 0.85416666666666663
$ws.Range("M2").Value = 0.865979381443299
$ws.Range("N2").Value = 0.87265917602996257
$ws.Range("O2").Value = 0.84615384615384615
$ws.Range("P2").Value = 0.89423076923076927
$ws.Range("Q2").Value = 0.87055016181229772
$ws.Range("R2").Value = 0.84838709677419355
$ws.Range("S2").Value = 0.86893203883495151

$ws.Range("L3").Value = 0.86813186813186816
$ws.Range("M3").Value = 0.82258064516129026
$ws.Range("N3").Value = 0.82417582417582413
$ws.Range("O3").Value = 0.85
$ws.Range("P3").Value = 0.82681564245810057
$ws.Range("Q3").Value = 0.83161512027491402
$ws.Range("R3").Value = 0.84563758389261745
$ws.Range("S3").Value = 0.86580086580086579

$ws.Range("A1:S3").Select() | Out-Null
